# Updated RAD EL-Motor Fuel Tax.
# Regenerated the "Date" column (B2:B48) of Sheet1 with a fresh batch of
# execution timestamps, as produced by a re-run of the Katalon RAD test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
    "Mon Apr 08 17:37:46 EDT 2024",
    "Mon Apr 08 17:37:55 EDT 2024",
    "Mon Apr 08 17:38:04 EDT 2024",
    "Mon Apr 08 17:38:13 EDT 2024",
    "Mon Apr 08 17:38:23 EDT 2024",
    "Mon Apr 08 17:38:32 EDT 2024",
    "Mon Apr 08 17:38:42 EDT 2024",
    "Mon Apr 08 17:38:51 EDT 2024",
    "Mon Apr 08 17:39:00 EDT 2024",
    "Mon Apr 08 17:39:10 EDT 2024",
    "Mon Apr 08 17:39:19 EDT 2024",
    "Mon Apr 08 17:39:28 EDT 2024",
    "Mon Apr 08 17:39:38 EDT 2024",
    "Mon Apr 08 17:39:47 EDT 2024",
    "Mon Apr 08 17:39:57 EDT 2024",
    "Mon Apr 08 17:40:06 EDT 2024",
    "Mon Apr 08 17:40:15 EDT 2024",
    "Mon Apr 08 17:40:25 EDT 2024",
    "Mon Apr 08 17:40:34 EDT 2024",
    "Mon Apr 08 17:40:43 EDT 2024",
    "Mon Apr 08 17:40:52 EDT 2024",
    "Mon Apr 08 17:41:03 EDT 2024",
    "Mon Apr 08 17:41:12 EDT 2024",
    "Mon Apr 08 17:41:21 EDT 2024",
    "Mon Apr 08 17:41:31 EDT 2024",
    "Mon Apr 08 17:41:41 EDT 2024",
    "Mon Apr 08 17:41:50 EDT 2024",
    "Mon Apr 08 17:42:00 EDT 2024",
    "Mon Apr 08 17:42:09 EDT 2024",
    "Mon Apr 08 17:42:19 EDT 2024",
    "Mon Apr 08 17:42:28 EDT 2024",
    "Mon Apr 08 17:42:37 EDT 2024",
    "Mon Apr 08 17:42:47 EDT 2024",
    "Mon Apr 08 17:42:56 EDT 2024",
    "Mon Apr 08 17:43:05 EDT 2024",
    "Mon Apr 08 17:43:15 EDT 2024",
    "Mon Apr 08 17:43:24 EDT 2024",
    "Mon Apr 08 17:43:34 EDT 2024",
    "Mon Apr 08 17:43:43 EDT 2024",
    "Mon Apr 08 17:43:52 EDT 2024",
    "Mon Apr 08 17:44:02 EDT 2024",
    "Mon Apr 08 17:44:11 EDT 2024",
    "Mon Apr 08 17:44:20 EDT 2024",
    "Mon Apr 08 17:44:30 EDT 2024",
    "Mon Apr 08 17:44:39 EDT 2024",
    "Mon Apr 08 17:44:49 EDT 2024",
    "Mon Apr 08 17:44:58 EDT 2024"
)

$startRow = 2
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
}
